# Auto-generated: update cached market-price / profit figures per commit
# "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H26").Value = 1000
$ws.Range("J26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("N26").Value = -1688

$ws.Range("H34").Value = 3157.4443
$ws.Range("I34").Value = 3157.4443
$ws.Range("K34").Value = 3157.4443
$ws.Range("M34").Value = -2954.4443

$ws.Range("H36").Value = 3157.4443
$ws.Range("I36").Value = 3157.4443
$ws.Range("K36").Value = 3157.4443
$ws.Range("M36").Value = -2442.4443

$ws.Range("H86").Value = 2314.7693
$ws.Range("I86").Value = 2244.4375
$ws.Range("J86").Value = 2427.3
$ws.Range("K86").Value = 2244.4375
$ws.Range("L86").Value = 2427.3
$ws.Range("M86").Value = -1121.4375
$ws.Range("N86").Value = -4673.3

$ws.Range("H89").Value = 2314.7693
$ws.Range("I89").Value = 2244.4375
$ws.Range("J89").Value = 2427.3
$ws.Range("K89").Value = 11222.1875
$ws.Range("L89").Value = 12136.5
$ws.Range("M89").Value = -5606.1875
$ws.Range("N89").Value = -23368.5

$ws.Range("H112").Value = 2879.7727
$ws.Range("I112").Value = 2026.1111
$ws.Range("J112").Value = 3470.7693
$ws.Range("K112").Value = 6078.3333
$ws.Range("L112").Value = 10412.3079
$ws.Range("M112").Value = -4970.3333
$ws.Range("N112").Value = -12628.3079

$ws.Range("H137").Value = 34488184
$ws.Range("J137").Value = 4012.8572
$ws.Range("L137").Value = 12038.5716
$ws.Range("N137").Value = -17138.5716


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 287179.12
$ws.Range("I32").Value = 461793.12
$ws.Range("K32").Value = 461793.12
$ws.Range("M32").Value = -461506.12

$ws.Range("H97").Value = 2102
$ws.Range("I97").Value = 2102
$ws.Range("K97").Value = 2102
$ws.Range("M97").Value = -1606

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H121").Value = 33994.5
$ws.Range("J121").Value = 59989
$ws.Range("L121").Value = 59989
$ws.Range("N121").Value = -63483


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 354.2857
$ws.Range("I12").Value = 256.25
$ws.Range("K12").Value = 256.25
$ws.Range("M12").Value = -88.25

$ws.Range("H94").Value = 2606.7
$ws.Range("I94").Value = 2396.3333
$ws.Range("J94").Value = 4500
$ws.Range("K94").Value = 2396.3333
$ws.Range("L94").Value = 4500
$ws.Range("M94").Value = -1945.3333
$ws.Range("N94").Value = -5402

$ws.Range("H128").Value = 7063.3335
$ws.Range("I128").Value = 7063.3335
$ws.Range("K128").Value = 21190.0005
$ws.Range("M128").Value = -18700.0005


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 91243.37
$ws.Range("J22").Value = 500750
$ws.Range("L22").Value = 500750
$ws.Range("N22").Value = -501450

$ws.Range("H31").Value = 1854475.8
$ws.Range("J31").Value = 3199.8
$ws.Range("L31").Value = 3199.8
$ws.Range("N31").Value = -3789.8

$ws.Range("H34").Value = 1854475.8
$ws.Range("J34").Value = 3199.8
$ws.Range("L34").Value = 3199.8
$ws.Range("N34").Value = -3603.8

$ws.Range("H132").Value = 2771.2632
$ws.Range("J132").Value = 3149
$ws.Range("L132").Value = 9447
$ws.Range("N132").Value = -14507


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 51
$ws.Range("I7").Value = 51
$ws.Range("K7").Value = 153
$ws.Range("M7").Value = -41

$ws.Range("H92").Value = 345.91666
$ws.Range("J92").Value = 347.36365
$ws.Range("L92").Value = 1042.09095
$ws.Range("N92").Value = -3538.09095

$ws.Range("H107").Value = 6100
$ws.Range("I107").Value = 666.6667
$ws.Range("J107").Value = 8428.571
$ws.Range("K107").Value = 2000.0001
$ws.Range("L107").Value = 25285.713
$ws.Range("M107").Value = -80.00009999999997
$ws.Range("N107").Value = -29125.713

$ws.Range("H122").Value = 1467323.9


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5999.875
$ws.Range("I12").Value = 4833.1665
$ws.Range("J12").Value = 9500
$ws.Range("K12").Value = 4833.1665
$ws.Range("L12").Value = 9500
$ws.Range("M12").Value = -4693.1665
$ws.Range("N12").Value = -9780

$ws.Range("H80").Value = 4148.3447
$ws.Range("I80").Value = 2943.1365
$ws.Range("J80").Value = 7936.143
$ws.Range("K80").Value = 2943.1365
$ws.Range("L80").Value = 7936.143
$ws.Range("M80").Value = -1945.1365
$ws.Range("N80").Value = -9932.143

$ws.Range("H83").Value = 4148.3447
$ws.Range("I83").Value = 2943.1365
$ws.Range("J83").Value = 7936.143
$ws.Range("K83").Value = 14715.6825
$ws.Range("L83").Value = 39680.715
$ws.Range("M83").Value = -9723.682500000001
$ws.Range("N83").Value = -49664.715


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 5979
$ws.Range("I23").Value = 5979
$ws.Range("K23").Value = 5979
$ws.Range("M23").Value = -5749

$ws.Range("H40").Value = 3782.6667
$ws.Range("I40").Value = 4014
$ws.Range("K40").Value = 4014
$ws.Range("M40").Value = -3878

$ws.Range("H93").Value = 3112.087
$ws.Range("I93").Value = 2728.7
$ws.Range("J93").Value = 5668
$ws.Range("K93").Value = 2728.7
$ws.Range("L93").Value = 5668
$ws.Range("M93").Value = -1480.7
$ws.Range("N93").Value = -8164


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H81").Value = 4767.2
$ws.Range("I81").Value = 4871.6665
$ws.Range("J81").Value = 4681.727
$ws.Range("K81").Value = 9743.333000000001
$ws.Range("L81").Value = 9363.454
$ws.Range("M81").Value = -8682.333000000001
$ws.Range("N81").Value = -11485.454

$ws.Range("H84").Value = 4767.2
$ws.Range("I84").Value = 4871.6665
$ws.Range("J84").Value = 4681.727
$ws.Range("K84").Value = 48716.665
$ws.Range("L84").Value = 46817.27
$ws.Range("M84").Value = -43412.665
$ws.Range("N84").Value = -57425.27

$ws.Range("H113").Value = 1340.1875
$ws.Range("I113").Value = 1564.875
$ws.Range("J113").Value = 1115.5
$ws.Range("K113").Value = 4694.625
$ws.Range("L113").Value = 3346.5
$ws.Range("M113").Value = -2524.625
$ws.Range("N113").Value = -7686.5

